# Apply the KeyItem.xlsx update: the data rows (A3:H16) were re-sorted by
# column C ("name") ascending instead of column B, the active selection
# moved to D9, and the row that now lands on row 16 (id 23, the wrapped
# "detail" cell) grows to a taller auto-fit height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:H16")
$sortKey   = $ws.Range("C2")

# Re-sort the table by column C (name) ascending via the sheet's AutoFilter,
# which keeps the sort metadata attached to the AutoFilter (matching how
# the workbook was actually re-sorted) instead of creating a second,
# unrelated sort-state entry on the worksheet.
$af = $ws.AutoFilter
$af.Sort.SortFields.Clear()
$af.Sort.SortFields.Add($sortKey, 0, 1, 0, 0)
$af.Sort.SetRange($dataRange)
$af.Sort.Header = 1
$af.Sort.Apply()

# Row heights are tied to the physical row, not the data, so after the sort
# the explicit wrapped-text height needs to be moved from its old location
# (row 6) to its new one (row 16).
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(16).RowHeight = 60

# Move the selection to where it ended up in the saved file.
[void]$ws.Range("D9").Select()
